$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = '37.048.74'
$ws.Range("E2").Value2 = '  -0.33%  '
$ws.Range("D3").Value2 = '2.004.22'
$ws.Range("E3").Value2 = '  -0.76%  '
$ws.Range("E4").Value2 = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = '257.41'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = '  +4.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = '0.618'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = '  -1.35%  '
$ws.Range("E7").Value2 = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = '55.83'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value2 = '  -7.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = '0.377'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value2 = '  -3.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = '0.0766'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = '  -5.06%  '
$ws.Range("E11").Value2 = '  -2.64%  '
$ws.Range("B12").Value2 = 'WrappedliquidstakedEther2.0'
$ws.Range("C12").Value2 = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D12").Value2 = '2.302.31'
$ws.Range("E12").Value2 = '  -0.67%  '
$ws.Range("B13").Value2 = 'Chainlink'
$ws.Range("C13").Value2 = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = '14.21'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = '  -5.16%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = '21.36'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = '  -2.32%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = '0.792'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = '  -6.53%  '
$ws.Range("B16").Value2 = 'Polkadot'
$ws.Range("C16").Value2 = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = '5.18'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = '  -4.66%  '
$ws.Range("B17").Value2 = 'WrappedEther'
$ws.Range("C17").Value2 = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value2 = '1.966.57'
$ws.Range("E17").Value2 = '  -2.73%  '
$ws.Range("D18").Value2 = '36.997.44'
$ws.Range("E18").Value2 = '  -0.40%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = '70.72'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = '  +0.60%  '
$ws.Range("D20").Value2 = '0.0₃0829'
$ws.Range("E20").Value2 = '  -3.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = '233.35'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = '  +1.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = '5.07'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = '  -2.78%  '
$ws.Range("E23").Value2 = '  +0.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = '2.55'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = '  -0.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = '2.36'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = '  +0.44%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = '164.45'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = '  +0.62%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = '8.89'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = '  -4.95%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = '19.45'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = '  -1.52%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = '1.34'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = '  -3.72%  '
$ws.Range("E30").Value2 = '  -9.68%  '
$ws.Range("E31").Value2 = '  -2.05%  '
$ws.Range("E32").Value2 = '  -4.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = '0.0637'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = '  -5.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = '4.40'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = '  -1.29%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = '2.34'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = '  -6.18%  '
$ws.Range("E36").Value2 = '  -4.44%  '
$ws.Range("E37").Value2 = '  +0.75%  '
$ws.Range("E38").Value2 = '  +0.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = '5.57'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = '  +4.66%  '
$ws.Range("E40").Value2 = '  +0.80%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = '1.17'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = '  -0.37%  '
$ws.Range("D42").Value2 = '1.438.11'
$ws.Range("E42").Value2 = '  +4.49%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = '0.0918'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = '  -5.47%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = '0.0209'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = '  -2.94%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = '88.79'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = '  -2.66%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = '15.49'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = '  -7.30%  '
$ws.Range("E47").Value2 = '  -3.31%  '
$ws.Range("E48").Value2 = '  +1.72%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = '6.91'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = '  -6.71%  '
$ws.Range("D50").Value2 = '2.194.61'
$ws.Range("E50").Value2 = '  -0.62%  '
$ws.Range("E51").Value2 = '  -8.39%  '
